$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "data refreshed" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 10:33"

# Row 7: Rusia - refreshed stats
$ws.Range("B7").Value = 771546
$ws.Range("C7").Value = 6109
$ws.Range("D7").Value = 550344
$ws.Range("E7").Value = 208860
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 95
$ws.Range("H7").Value = 12342

# Rows 33-35: Filipinas moves up above Bielorrusia/Oman with refreshed stats;
# Bielorrusia and Oman keep their previous values but shift down a row.
$ws.Range("A33").Value = "Filipinas"
$ws.Range("B33").Value = 67456
$ws.Range("C33").Value = 2241
$ws.Range("D33").Value = 22465
$ws.Range("E33").Value = 43160
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 58
$ws.Range("H33").Value = 1831

$ws.Range("A34").Value = "Bielorrusia"
$ws.Range("B34").Value = 65953
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 57856
$ws.Range("E34").Value = 7602
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 495

$ws.Range("A35").Value = "Oman"
$ws.Range("B35").Value = 65504
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 42772
$ws.Range("E35").Value = 22424
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 308

# Row 44: Israel - refreshed stats
$ws.Range("B44").Value = 49575
$ws.Range("C44").Value = 210
$ws.Range("D44").Value = 21440
$ws.Range("E44").Value = 27729
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 5
$ws.Range("H44").Value = 406

# Row 46: Singapur - refreshed stats
$ws.Range("B46").Value = 47912
$ws.Range("C46").Value = 257
$ws.Range("E46").Value = 4052

# Row 47: Polonia - refreshed stats
$ws.Range("D47").Value = 30292
$ws.Range("E47").Value = 7836

# Row 121: Eslovaquia - refreshed stats
$ws.Range("B121").Value = 1979
$ws.Range("C121").Value = 3
$ws.Range("D121").Value = 1530
$ws.Range("E121").Value = 421
